$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 2.3
$ws.Range("I4").Value = 2.67
$ws.Range("J4").Value = 2.72
$ws.Range("L4").Value = 3.05
$ws.Range("R4").Value = 2.67
$ws.Range("Z4").Value = 27
$ws.Range("AA4").Value = 15.5
$ws.Range("AB4").Value = 17.5
$ws.Range("AH4").Value = 15.5
$ws.Range("AI4").Value = 19.5
$ws.Range("AJ4").Value = 10.25
$ws.Range("AL4").Value = 18.5
$ws.Range("AM4").Value = 19
$ws.Range("AN4").Value = 4.9
$ws.Range("AO4").Value = 11.5
$ws.Range("AP4").Value = 14
$ws.Range("AQ4").Value = 37
$ws.Range("AS4").Value = 110
$ws.Range("AW4").Value = 5.3
$ws.Range("AX4").Value = 13.5
$ws.Range("AY4").Value = 15
$ws.Range("AZ4").Value = 50
$ws.Range("BA4").Value = 55

# Row 6 updates
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 8.5
$ws.Range("Q6").Value = 1.83
$ws.Range("R6").Value = 1.98
